# Apply updated Betfair back/lay odds for 2025-11-18 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "09:00:00"
$ws.Range("F2").Value = 2.74
$ws.Range("G2").Value = 3.2
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 2.68
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1.53
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 2.84
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 1.66
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.23
$ws.Range("S2").Value = 4.7
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 1.81
$ws.Range("V2").Value = 1.41
$ws.Range("W2").Value = 1.45
$ws.Range("AB2").Value = 12
$ws.Range("AC2").Value = 7.2
$ws.Range("AH2").Value = 60

# Row 4
$ws.Range("G4").Value = 2.28
$ws.Range("H4").Value = 3.2
$ws.Range("L4").Value = 1.34
$ws.Range("N4").Value = 4.6
$ws.Range("P4").Value = 2.24
$ws.Range("Q4").Value = 1.72
$ws.Range("R4").Value = 1.49
$ws.Range("S4").Value = 2.78
$ws.Range("T4").Value = 1.62
$ws.Range("U4").Value = 2.36
$ws.Range("W4").Value = 1.78
$ws.Range("AD4").Value = 14.5
$ws.Range("AG4").Value = 11.5
$ws.Range("AH4").Value = 16.5
$ws.Range("AI4").Value = 100
$ws.Range("AM4").Value = 580
$ws.Range("AN4").Value = 13
$ws.Range("AO4").Value = 200

# Row 5
$ws.Range("F5").Value = 2.06
$ws.Range("G5").Value = 2.18
$ws.Range("H5").Value = 3.95
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.55
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 3.25
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 1.76
$ws.Range("R5").Value = 1.29
$ws.Range("S5").Value = 3.8
$ws.Range("V5").Value = 1.28
$ws.Range("W5").Value = 1.84
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 30
$ws.Range("AA5").Value = 200
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 17.5
$ws.Range("AF5").Value = 13
$ws.Range("AJ5").Value = 29
$ws.Range("AK5").Value = 70
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 55
$ws.Range("AO5").Value = 75

# Row 6
$ws.Range("F6").Value = 1.25
$ws.Range("G6").Value = 1.28
$ws.Range("H6").Value = 15.5
$ws.Range("I6").Value = 17.5
$ws.Range("J6").Value = 6.6
$ws.Range("K6").Value = 7
$ws.Range("L6").Value = 1.31
$ws.Range("N6").Value = 5.3
$ws.Range("O6").Value = 1.21
$ws.Range("P6").Value = 2.44
$ws.Range("Q6").Value = 1.65
$ws.Range("R6").Value = 1.58
$ws.Range("S6").Value = 2.64
$ws.Range("T6").Value = 2.24
$ws.Range("U6").Value = 1.71
$ws.Range("W6").Value = 4.6
$ws.Range("X6").Value = 25
$ws.Range("Y6").Value = 50
$ws.Range("Z6").Value = 170
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 9
$ws.Range("AC6").Value = 15
$ws.Range("AD6").Value = 60
$ws.Range("AE6").Value = 320
$ws.Range("AF6").Value = 7.4
$ws.Range("AG6").Value = 11
$ws.Range("AI6").Value = 230
$ws.Range("AJ6").Value = 9
$ws.Range("AK6").Value = 14
$ws.Range("AL6").Value = 44
$ws.Range("AM6").Value = 260
$ws.Range("AN6").Value = 4.5
$ws.Range("AO6").Value = 480

# Row 7
$ws.Range("F7").Value = 2.08
$ws.Range("G7").Value = 2.22
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 4.5
$ws.Range("K7").Value = 3.55
$ws.Range("L7").Value = 1.5
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 2.96
$ws.Range("P7").Value = 1.68
$ws.Range("Q7").Value = 2.26
$ws.Range("S7").Value = 4.3
$ws.Range("T7").Value = 1.95
$ws.Range("V7").Value = 1.29
$ws.Range("W7").Value = 1.81
$ws.Range("X7").Value = 11
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 130
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 18
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 65
$ws.Range("AN7").Value = 24
$ws.Range("AO7").Value = 100

# Row 8
$ws.Range("F8").Value = 2.26
$ws.Range("G8").Value = 2.34
$ws.Range("H8").Value = 3.75
$ws.Range("K8").Value = 3.25
$ws.Range("L8").Value = 1.51
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 3
$ws.Range("O8").Value = 1.45
$ws.Range("Q8").Value = 2.36
$ws.Range("R8").Value = 1.24
$ws.Range("S8").Value = 4.7
$ws.Range("T8").Value = 2.02
$ws.Range("U8").Value = 1.85
$ws.Range("W8").Value = 1.74
$ws.Range("X8").Value = 10.5
$ws.Range("Z8").Value = 32
$ws.Range("AB8").Value = 8.6
$ws.Range("AD8").Value = 17.5
$ws.Range("AE8").Value = 65

# Row 9
$ws.Range("F9").Value = 1.19
$ws.Range("H9").Value = 21
$ws.Range("I9").Value = 26
$ws.Range("K9").Value = 8.800000000000001
$ws.Range("L9").Value = 1.3
$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 2.46
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 1.57
$ws.Range("S9").Value = 2.48
$ws.Range("V9").Value = 1.04
$ws.Range("X9").Value = 30
$ws.Range("Y9").Value = 70
$ws.Range("Z9").Value = 310
$ws.Range("AB9").Value = 9
$ws.Range("AC9").Value = 19.5
$ws.Range("AD9").Value = 95
$ws.Range("AF9").Value = 7.8
$ws.Range("AJ9").Value = 8
$ws.Range("AL9").Value = 70
$ws.Range("AM9").Value = 490
$ws.Range("AN9").Value = 4.1

# Row 10
$ws.Range("G10").Value = 1.31
$ws.Range("J10").Value = 6
$ws.Range("N10").Value = 4.5
$ws.Range("O10").Value = 1.23
$ws.Range("P10").Value = 2.2
$ws.Range("Q10").Value = 1.73
$ws.Range("R10").Value = 1.46
$ws.Range("S10").Value = 2.8
$ws.Range("T10").Value = 2.3
$ws.Range("U10").Value = 1.64
$ws.Range("V10").Value = 1.06
$ws.Range("Z10").Value = 180
$ws.Range("AC10").Value = 15
$ws.Range("AE10").Value = 370
$ws.Range("AI10").Value = 270
$ws.Range("AK10").Value = 15.5
$ws.Range("AM10").Value = 320
$ws.Range("AN10").Value = 5.3

# Row 11
$ws.Range("F11").Value = 1.85
$ws.Range("G11").Value = 1.91
$ws.Range("I11").Value = 5.1
$ws.Range("J11").Value = 3.65
$ws.Range("K11").Value = 3.8
$ws.Range("N11").Value = 3.3
$ws.Range("P11").Value = 1.79
$ws.Range("S11").Value = 3.9
$ws.Range("T11").Value = 1.95
$ws.Range("U11").Value = 1.88
$ws.Range("V11").Value = 1.24
$ws.Range("W11").Value = 2.08
$ws.Range("X11").Value = 13
$ws.Range("Y11").Value = 15
$ws.Range("Z11").Value = 44
$ws.Range("AA11").Value = 140
$ws.Range("AB11").Value = 8
$ws.Range("AC11").Value = 8.6
$ws.Range("AE11").Value = 100
$ws.Range("AF11").Value = 11
$ws.Range("AI11").Value = 90
$ws.Range("AJ11").Value = 21
$ws.Range("AL11").Value = 980
$ws.Range("AM11").Value = 170
$ws.Range("AN11").Value = 15.5
$ws.Range("AO11").Value = 110

# Row 12
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 2.04
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 3.35
$ws.Range("K12").Value = 3.5
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 3.05
$ws.Range("P12").Value = 1.69
$ws.Range("Q12").Value = 2.26
$ws.Range("T12").Value = 1.98
$ws.Range("U12").Value = 1.83
$ws.Range("V12").Value = 1.27
$ws.Range("W12").Value = 1.96
$ws.Range("X12").Value = 11.5
$ws.Range("AA12").Value = 130
$ws.Range("AC12").Value = 7.8
$ws.Range("AD12").Value = 19.5
$ws.Range("AF12").Value = 11.5
$ws.Range("AN12").Value = 20
$ws.Range("AO12").Value = 95

# Row 13
$ws.Range("F13").Value = 1.09
$ws.Range("G13").Value = 1.13
$ws.Range("H13").Value = 27
$ws.Range("I13").Value = 38
$ws.Range("J13").Value = 11.5
$ws.Range("K13").Value = 18
$ws.Range("L13").Value = 1.2
$ws.Range("N13").Value = 7.4
$ws.Range("O13").Value = 1.11
$ws.Range("P13").Value = 3.35
$ws.Range("Q13").Value = 1.33
$ws.Range("R13").Value = 2.02
$ws.Range("S13").Value = 1.86
$ws.Range("T13").Value = 2.4
$ws.Range("U13").Value = 1.54
$ws.Range("W13").Value = 8.4
$ws.Range("Y13").Value = 120
$ws.Range("Z13").Value = 400
$ws.Range("AB13").Value = 13.5
$ws.Range("AC13").Value = 42
$ws.Range("AD13").Value = 130
$ws.Range("AF13").Value = 9.199999999999999
$ws.Range("AG13").Value = 17
$ws.Range("AJ13").Value = 8.6
$ws.Range("AK13").Value = 16.5
$ws.Range("AM13").Value = 400
$ws.Range("AN13").Value = 2.8

# Row 14
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 4.3
$ws.Range("H14").Value = 2.06
$ws.Range("I14").Value = 2.12
$ws.Range("K14").Value = 3.75
$ws.Range("P14").Value = 1.6
$ws.Range("Q14").Value = 2.42
